$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 684
$ws.Range("I19").Value = 456.5
$ws.Range("J19").Value = 957
$ws.Range("K19").Value = 456.5
$ws.Range("L19").Value = 957
$ws.Range("M19").Value = -281.5
$ws.Range("N19").Value = -1307

# Row 96
$ws.Range("H96").Value = 2535.3684
$ws.Range("I96").Value = 2321
$ws.Range("K96").Value = 6963
$ws.Range("M96").Value = -5590

# Row 111
$ws.Range("H111").Value = 718.3333
$ws.Range("I111").Value = 577.5
$ws.Range("K111").Value = 1732.5
$ws.Range("M111").Value = 1334.5

# Row 113
$ws.Range("H113").Value = 3100.4167
$ws.Range("I113").Value = 3400.625
$ws.Range("J113").Value = 2500
$ws.Range("K113").Value = 3400.625
$ws.Range("L113").Value = 2500
$ws.Range("M113").Value = -146.625
$ws.Range("N113").Value = -9008

# Row 114
$ws.Range("H114").Value = 25087.143
$ws.Range("J114").Value = 25087.143
$ws.Range("L114").Value = 25087.143
$ws.Range("N114").Value = -33765.143

# Row 116
$ws.Range("H116").Value = 2456.1667
$ws.Range("I116").Value = 2362.6924
$ws.Range("J116").Value = 2699.2
$ws.Range("K116").Value = 2362.6924
$ws.Range("L116").Value = 2699.2
$ws.Range("M116").Value = 1079.3076
$ws.Range("N116").Value = -9583.200000000001

$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 718.725
$ws.Range("I74").Value = 746.3
$ws.Range("J74").Value = 691.15
$ws.Range("K74").Value = 746.3
$ws.Range("L74").Value = 691.15
$ws.Range("M74").Value = 127.7
$ws.Range("N74").Value = -2439.15

# Row 77
$ws.Range("H77").Value = 718.725
$ws.Range("I77").Value = 746.3
$ws.Range("J77").Value = 691.15
$ws.Range("K77").Value = 3731.5
$ws.Range("L77").Value = 3455.75
$ws.Range("M77").Value = 636.5
$ws.Range("N77").Value = -12191.75

# Row 97
$ws.Range("H97").Value = 327.56522
$ws.Range("I97").Value = 316.2
$ws.Range("J97").Value = 403.33334
$ws.Range("K97").Value = 316.2
$ws.Range("L97").Value = 403.33334
$ws.Range("M97").Value = 179.8
$ws.Range("N97").Value = -1395.33334

# Row 102
$ws.Range("H102").Value = 3269.1667
$ws.Range("I102").Value = 2866.25
$ws.Range("J102").Value = 4075
$ws.Range("K102").Value = 2866.25
$ws.Range("L102").Value = 4075
$ws.Range("M102").Value = -1244.25
$ws.Range("N102").Value = -7319

# Row 110
$ws.Range("H110").Value = 639.75
$ws.Range("J110").Value = 399.5
$ws.Range("L110").Value = 399.5
$ws.Range("N110").Value = -4489.5

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 766.1539
$ws.Range("I94").Value = 442.85715
$ws.Range("K94").Value = 442.85715
$ws.Range("M94").Value = 8.14285000000001

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 233.93333
$ws.Range("I22").Value = 223.76923
$ws.Range("J22").Value = 300
$ws.Range("K22").Value = 223.76923
$ws.Range("L22").Value = 300
$ws.Range("M22").Value = 126.23077
$ws.Range("N22").Value = -1000

# Row 122
$ws.Range("H122").Value = 914038.75
$ws.Range("J122").Value = 12399.75
$ws.Range("L122").Value = 37199.25
$ws.Range("N122").Value = -42099.25

$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 4940.3477
$ws.Range("J3").Value = 8366.182000000001
$ws.Range("L3").Value = 25098.546
$ws.Range("N3").Value = -25322.546

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 1580.3
$ws.Range("I7").Value = 1484
$ws.Range("J7").Value = 1724.75
$ws.Range("K7").Value = 1484
$ws.Range("L7").Value = 1724.75
$ws.Range("M7").Value = -1372
$ws.Range("N7").Value = -1948.75

# Row 40
$ws.Range("H40").Value = 2325.5
$ws.Range("I40").Value = 2329.1428
$ws.Range("K40").Value = 2329.1428
$ws.Range("M40").Value = -2193.1428

# Row 46
$ws.Range("J46").Value = 500
$ws.Range("L46").Value = 500
$ws.Range("N46").Value = -876

# Row 61
$ws.Range("H61").Value = 17633
$ws.Range("I61").Value = 17633
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 17633
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -17431
$ws.Range("N61").ClearContents()

# Row 113
$ws.Range("H113").Value = 17633
$ws.Range("I113").Value = 17633
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 17633
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -15463
$ws.Range("N113").ClearContents()

# Row 126
$ws.Range("H126").Value = 1580.3
$ws.Range("I126").Value = 1484
$ws.Range("J126").Value = 1724.75
$ws.Range("K126").Value = 4452
$ws.Range("L126").Value = 5174.25
$ws.Range("M126").Value = -1982
$ws.Range("N126").Value = -10114.25

$ws = $wb.Worksheets.Item("WVR")
# Row 39
$ws.Range("H39").Value = 6044
$ws.Range("I39").Value = 6044
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 6044
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -5631
$ws.Range("N39").ClearContents()

# Row 42
$ws.Range("H42").Value = 33361650
$ws.Range("I42").Value = 100000050
$ws.Range("J42").Value = 42450
$ws.Range("K42").Value = 100000050
$ws.Range("L42").Value = 42450
$ws.Range("M42").Value = -99999672
$ws.Range("N42").Value = -43206

# Row 64
$ws.Range("H64").Value = 19000
$ws.Range("J64").Value = 19000
$ws.Range("L64").Value = 19000
$ws.Range("N64").Value = -19496

# Row 67
$ws.Range("H67").Value = 19000
$ws.Range("J67").Value = 19000
$ws.Range("L67").Value = 19000
$ws.Range("N67").Value = -20716

# Row 70
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

# Row 73
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

# Row 81
$ws.Range("I81").Value = 1236.5714
$ws.Range("J81").Value = 1828.8572
$ws.Range("K81").Value = 2473.1428
$ws.Range("L81").Value = 3657.7144
$ws.Range("M81").Value = -1412.1428
$ws.Range("N81").Value = -5779.7144

# Row 84
$ws.Range("I84").Value = 1236.5714
$ws.Range("J84").Value = 1828.8572
$ws.Range("K84").Value = 12365.714
$ws.Range("L84").Value = 18288.572
$ws.Range("M84").Value = -7061.714
$ws.Range("N84").Value = -28896.572

# Row 122
$ws.Range("H122").Value = 1499.238
$ws.Range("I122").Value = 1232.2667
$ws.Range("J122").Value = 2166.6667
$ws.Range("K122").Value = 3696.800099999999
$ws.Range("L122").Value = 6500.000100000001
$ws.Range("M122").Value = -1246.800099999999
$ws.Range("N122").Value = -11400.0001
